$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.799.04'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '1.702.13'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '317.40'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3957'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4087'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.73%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.507'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.004'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '52.68'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08907'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.708'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +6.76%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '24.35'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.54%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.163'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001334'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').Value = '1.707.02'
$ws.Range('E17').Value = '  +0.86%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '99.88'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07138'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '20.07'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.218'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.07%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.007'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.73%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.65'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.30%  '
$ws.Range('D24').Value = '24.785.12'
$ws.Range('E24').Value = '  +0.72%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.104'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.05%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.338'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.99'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.91%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.334'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +24.68%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '165.11'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '139.72'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.25%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.198'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.062'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +12.70%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09204'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +7.16%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.084'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.55%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.03051'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +11.25%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.2814'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.67%  '
$ws.Range('E37').Value = '  +1.81%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '11.07'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.96%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '14.60'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.09309'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.93%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.7838'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.16%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.479'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '16.26'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.55%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.638'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.38%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.7264'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.249'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.361'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.62%  '
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '141.59'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '93.36'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +4.96%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.08074'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.04%  '
